$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four task-description cells: ".class" -> ".java"
$ws.Range("C5").Value = "Interface utilisateur (Application.java, WindowPrincipal.java, WindowMenu)"
$ws.Range("J5").Value = "Interface utilisteur (WindowMenu.java, WindowTable.java, WindowButton)"
$ws.Range("J11").Value = "Implémentation des fonctionalités(OpenAndSave.java)"
$ws.Range("B11").Value = "Implémentation des fonctionalités (WindowPrincipal.java,  WindowTable.java)"

# Update the sheet view: zoom level and active cell selection
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("H15").Select() | Out-Null
